$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = $false
$ws.Range("B1").Value = "simple"
$ws.Range("C1").Value = 2995.996492127714

$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "simple"
$ws.Range("C2").Value = 2995.996492127714
